# Fix the 2050 column label (was showing a stray numeric placeholder) and
# remove the "Total" summary rows from the tables.

$wb = $excel.ActiveWorkbook

# Sheets 1-3 and 5 ("...2050" column) should read "2050" as text instead of
# the stray leftover numeric placeholder (671.09...).
# Sheet 4 uses year-range headers ("2015-2030", "2031-2040", ...), so its
# last column should read "2041-2050" instead.
$yearSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

# -4122 == xlPasteFormats: re-stamp E1 with D1's exact formatting (bold,
# centered, bordered text style) after the value write so the cell keeps
# sharing the original style record instead of picking up a fresh one.
$xlPasteFormats = -4122

foreach ($name in $yearSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E1").NumberFormat = "@"
    $ws.Range("E1").Value = "2050"
    $ws.Range("D1").Copy()
    $ws.Range("E1").PasteSpecial($xlPasteFormats)
}

$wsIncremental = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIncremental.Range("E1").NumberFormat = "@"
$wsIncremental.Range("E1").Value = "2041-2050"
$wsIncremental.Range("D1").Copy()
$wsIncremental.Range("E1").PasteSpecial($xlPasteFormats)

# Remove the "Total" row from the tables that have one (row 13 for the
# 4 main tables, row 4 for the cost table). The emissions table never had
# a Total row, so it only needs the label fix above.
$rowsToDelete = @{
    "Potencia Acumulada - SIN (MW)"  = 13
    "Geracao Periodo Medio (MWMed)"  = 13
    "Atendimento a Ponta(MW)"        = 13
    "Potencia Incremental - SIN(MW)" = 13
    "Custo Total (bilhões de R$)"    = 4
}

foreach ($name in $rowsToDelete.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $rowNum = $rowsToDelete[$name]
    $ws.Rows.Item($rowNum).Delete()
}
